$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.474.96'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '1.906.43'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = "'0.704"
$ws.Range('E5').Value = '  +11.21%  '
$ws.Range('D6').Value = "'246.85"
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = "'40.78"
$ws.Range('E8').Value = '  -2.87%  '
$ws.Range('D9').Value = "'0.353"
$ws.Range('E9').Value = '  +3.86%  '
$ws.Range('D10').Value = "'52.62"
$ws.Range('E10').Value = '  +8.10%  '
$ws.Range('D11').Value = "'0.0727"
$ws.Range('E11').Value = '  +3.38%  '
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('D13').Value = '2.181.04'
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('E14').Value = '  +2.66%  '
$ws.Range('D15').Value = "'0.714"
$ws.Range('E15').Value = '  +2.97%  '
$ws.Range('D16').Value = '1.906.49'
$ws.Range('E16').Value = '  +0.65%  '
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('D18').Value = '35.427.59'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('D19').Value = "'73.07"
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('D20').Value = '0.0₃0825'
$ws.Range('E20').Value = '  +0.26%  '
$ws.Range('D21').Value = "'242.11"
$ws.Range('E21').Value = '  -0.43%  '
$ws.Range('D22').Value = "'12.91"
$ws.Range('E22').Value = '  +2.50%  '
$ws.Range('D23').Value = "'5.07"
$ws.Range('E23').Value = '  +4.38%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('E25').Value = '  +0.95%  '
$ws.Range('E26').Value = '  +5.69%  '
$ws.Range('D27').Value = "'169.12"
$ws.Range('E27').Value = '  -2.02%  '
$ws.Range('D28').Value = "'8.69"
$ws.Range('E28').Value = '  +1.66%  '
$ws.Range('D29').Value = "'18.95"
$ws.Range('E29').Value = '  +5.64%  '
$ws.Range('D30').Value = "'0.132"
$ws.Range('E30').Value = '  +5.27%  '
$ws.Range('E32').Value = '  +3.08%  '
$ws.Range('E33').Value = '  +0.63%  '
$ws.Range('E34').Value = '  +7.43%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = "'4.19"
$ws.Range('E35').Value = '  +0.95%  '
$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D36').Value = "'1.01"
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('D37').Value = "'0.921"
$ws.Range('E37').Value = '  -5.45%  '
$ws.Range('D38').Value = "'1.47"
$ws.Range('E38').Value = '  +9.11%  '
$ws.Range('D39').Value = "'2.02"
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('D40').Value = "'96.59"
$ws.Range('E40').Value = '  +6.43%  '
$ws.Range('D41').Value = "'1.12"
$ws.Range('E41').Value = '  +0.84%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = "'0.0655"
$ws.Range('E42').Value = '  +3.77%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').Value = "'16.52"
$ws.Range('E43').Value = '  +5.29%  '
$ws.Range('E44').Value = '  +1.67%  '
$ws.Range('D45').Value = '1.357.96'
$ws.Range('E45').Value = '  +0.67%  '
$ws.Range('D46').Value = "'2.42"
$ws.Range('E46').Value = '  +2.02%  '
$ws.Range('B47').Value = 'HuobiToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D47').Value = "'2.42"
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('B48').Value = 'MXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D48').Value = "'2.79"
$ws.Range('E48').Value = '  +1.34%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').Value = "'46.06"
$ws.Range('E49').Value = '  -8.90%  '
$ws.Range('D50').Value = "'12.21"
$ws.Range('E50').Value = '  -4.52%  '
$ws.Range('E51').Value = '  -1.30%  '
